$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "36"
$a8.Characters(21, 2).Font.Size = 10

$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "9/4/2023"
$c9.Characters(27, 8).Font.Size = 10
$c9.Characters(46, 8).Text = "9/10/2023"
$c9.Characters(46, 9).Font.Size = 10

# --- Crime statistics table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 13
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = 44.444444444444
$ws.Range("I14").Value = 92
$ws.Range("J14").Value = 99
$ws.Range("K14").Value = -7.070707070707
$ws.Range("L14").Value = -17.857142857142
$ws.Range("M14").Value = -4.166666666666
$ws.Range("N14").Value = -73.255813953488

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 34
$ws.Range("G15").Value = 24
$ws.Range("H15").Value = 41.666666666666
$ws.Range("I15").Value = 268
$ws.Range("J15").Value = 276
$ws.Range("K15").Value = -2.898550724637
$ws.Range("L15").Value = 8.502024291497
$ws.Range("M15").Value = 24.651162790697
$ws.Range("N15").Value = -46.930693069306

# Row 16
$ws.Range("C16").Value = 98
$ws.Range("D16").Value = 99
$ws.Range("E16").Value = -1.010101010101
$ws.Range("F16").Value = 422
$ws.Range("G16").Value = 406
$ws.Range("H16").Value = 3.940886699507
$ws.Range("I16").Value = 3368
$ws.Range("J16").Value = 3563
$ws.Range("K16").Value = -5.472916081953
$ws.Range("L16").Value = 34.129828753484
$ws.Range("M16").Value = 12.906469996647
$ws.Range("N16").Value = -69.608373939722

# Row 17
$ws.Range("C17").Value = 165
$ws.Range("D17").Value = 145
$ws.Range("E17").Value = 13.793103448275
$ws.Range("F17").Value = 671
$ws.Range("G17").Value = 580
$ws.Range("H17").Value = 15.689655172413
$ws.Range("I17").Value = 5643
$ws.Range("J17").Value = 5133
$ws.Range("K17").Value = 9.935710111046
$ws.Range("L17").Value = 33.404255319148
$ws.Range("M17").Value = 80.576
$ws.Range("N17").Value = -12.606473594548

# Row 18
$ws.Range("C18").Value = 48
$ws.Range("D18").Value = 50
$ws.Range("E18").Value = -4
$ws.Range("F18").Value = 214
$ws.Range("G18").Value = 217
$ws.Range("H18").Value = -1.382488479262
$ws.Range("I18").Value = 2062
$ws.Range("J18").Value = 2024
$ws.Range("K18").Value = 1.877470355731
$ws.Range("L18").Value = 40.846994535519
$ws.Range("M18").Value = -7.616487455197
$ws.Range("N18").Value = -84.185903826980

# Row 19
$ws.Range("C19").Value = 181
$ws.Range("D19").Value = 153
$ws.Range("E19").Value = 18.300653594771
$ws.Range("F19").Value = 699
$ws.Range("G19").Value = 631
$ws.Range("H19").Value = 10.776545166402
$ws.Range("I19").Value = 5499
$ws.Range("J19").Value = 5531
$ws.Range("K19").Value = -0.578557222925
$ws.Range("L19").Value = 25.376196990424
$ws.Range("M19").Value = 71.522145976294
$ws.Range("N19").Value = 4.523854780459

# Row 20
$ws.Range("C20").Value = 102
$ws.Range("D20").Value = 63
$ws.Range("E20").Value = 61.904761904761
$ws.Range("F20").Value = 404
$ws.Range("G20").Value = 253
$ws.Range("H20").Value = 59.683794466403
$ws.Range("I20").Value = 3710
$ws.Range("J20").Value = 2694
$ws.Range("K20").Value = 37.713437268003
$ws.Range("L20").Value = 92.927717108684
$ws.Range("M20").Value = 157.817929117443
$ws.Range("N20").Value = -64.970257766027

# Row 21
$ws.Range("C21").Value = 601
$ws.Range("D21").Value = 517
$ws.Range("E21").Value = 16.247582205029
$ws.Range("F21").Value = 2457
$ws.Range("G21").Value = 2120
$ws.Range("H21").Value = 15.896226415094
$ws.Range("I21").Value = 20642
$ws.Range("J21").Value = 19320
$ws.Range("K21").Value = 6.842650103519
$ws.Range("L21").Value = 38.788408525516
$ws.Range("M21").Value = 55.249699157641
$ws.Range("N21").Value = -56.340024112185

# Row 22
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = -22.222222222222
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 35
$ws.Range("H22").Value = -45.714285714285
$ws.Range("I22").Value = 199
$ws.Range("J22").Value = 253
$ws.Range("K22").Value = -21.343873517786
$ws.Range("L22").Value = 19.879518072289
$ws.Range("M22").Value = -8.294930875576

# Row 23
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 39
$ws.Range("E23").Value = -38.461538461538
$ws.Range("F23").Value = 131
$ws.Range("G23").Value = 123
$ws.Range("H23").Value = 6.504065040650
$ws.Range("I23").Value = 1239
$ws.Range("J23").Value = 1119
$ws.Range("K23").Value = 10.723860589812
$ws.Range("L23").Value = 51.466992665036
$ws.Range("M23").Value = 65.863453815261

# Row 24
$ws.Range("C24").Value = 308
$ws.Range("D24").Value = 324
$ws.Range("E24").Value = -4.938271604938
$ws.Range("F24").Value = 1395
$ws.Range("G24").Value = 1489
$ws.Range("H24").Value = -6.312961719274
$ws.Range("I24").Value = 12421
$ws.Range("J24").Value = 12838
$ws.Range("K24").Value = -3.248169496806
$ws.Range("L24").Value = 42.491682918435
$ws.Range("M24").Value = 40.540846345327

# Row 25
$ws.Range("C25").Value = 203
$ws.Range("D25").Value = 184
$ws.Range("E25").Value = 10.326086956521
$ws.Range("F25").Value = 836
$ws.Range("G25").Value = 717
$ws.Range("H25").Value = 16.596931659693
$ws.Range("I25").Value = 7367
$ws.Range("J25").Value = 6985
$ws.Range("K25").Value = 5.468861846814
$ws.Range("L25").Value = 26.320301783264
$ws.Range("M25").Value = -5.211013896037

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 55
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = 27.906976744186
$ws.Range("I26").Value = 453
$ws.Range("J26").Value = 478
$ws.Range("K26").Value = -5.230125523012
$ws.Range("L26").Value = 11.302211302211

# Row 27
$ws.Range("C27").Value = 27
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = 107.692307692308
$ws.Range("F27").Value = 85
$ws.Range("G27").Value = 71
$ws.Range("H27").Value = 19.718309859154
$ws.Range("I27").Value = 726
$ws.Range("J27").Value = 628
$ws.Range("K27").Value = 15.605095541401
$ws.Range("L27").Value = 15.055467511885

# Row 28
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 9
$ws.Range("E28").Value = -11.111111111111
$ws.Range("F28").Value = 38
$ws.Range("G28").Value = 30
$ws.Range("H28").Value = 26.666666666666
$ws.Range("I28").Value = 287
$ws.Range("J28").Value = 361
$ws.Range("K28").Value = -20.498614958448
$ws.Range("L28").Value = -33.870967741935
$ws.Range("M28").Value = -16.811594202898
$ws.Range("N28").Value = -70.921985815602

# Row 29
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 60
$ws.Range("F29").Value = 34
$ws.Range("G29").Value = 22
$ws.Range("H29").Value = 54.545454545454
$ws.Range("I29").Value = 237
$ws.Range("J29").Value = 305
$ws.Range("K29").Value = -22.295081967213
$ws.Range("L29").Value = -35.068493150684
$ws.Range("M29").Value = -17.421602787456
$ws.Range("N29").Value = -73.549107142857

# Row 30
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 15
$ws.Range("K30").Value = -55.882352941176
$ws.Range("L30").Value = -55.882352941176
